$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the Table-of-Contents entry "4. Feature Specifications" is
# styled "List Number" (there is a second, unrelated "Heading 1"
# paragraph with the same text further down that must NOT be touched).
# Append a manual line break followed by a bullet line announcing the
# new Factorial feature, matching:
#   <w:r>
#     <w:t>4. Feature Specifications</w:t>
#     <w:br/>
#     <w:t>\u2022 Factorial Operation: Calculates the factorial of a number</w:t>
#   </w:r>
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $bodyText = $p.Range.Text.TrimEnd([char]13)
    if ($p.Style.NameLocal -eq "List Number" -and $bodyText -eq "4. Feature Specifications") {
        $r = $p.Range
        $r.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
        $r.Collapse(0) | Out-Null      # wdCollapseEnd

        $bullet = [char]0x2022
        # "" + forces string concatenation (rather than numeric char addition)
        $newLine = "" + [char]11 + $bullet + " Factorial Operation: Calculates the factorial of a number"
        $r.InsertAfter($newLine)
        break
    }
}

# ---------------------------------------------------------------------
# Change 2: add a new "Factorial" / "factorial" row at the end of the
# two-column Menu Options table (the one whose header row is
# "Option" | "Function" and whose last existing row is "8" | "Exit").
# ---------------------------------------------------------------------
foreach ($tbl in $d.Tables) {
    $headerCell1 = $tbl.Rows.Item(1).Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7)
    $headerCell2 = $tbl.Rows.Item(1).Cells.Item(2).Range.Text.TrimEnd([char]13, [char]7)
    if ($headerCell1 -eq "Option" -and $headerCell2 -eq "Function") {
        $newRow = $tbl.Rows.Add()
        $newRow.Cells.Item(1).Range.Text = "Factorial"
        $newRow.Cells.Item(2).Range.Text = "factorial"
        break
    }
}
